$d = $word.ActiveDocument

# 1. Replace the placeholder ID text in the first paragraph and absorb the
#    trailing space-only run that followed it, so that run disappears.
$null = $d.Content.Find.Execute("**ID__AFFARS_5349_topic_7__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5349_402_6__ID**", 2)

# 2. Update paragraph formatting for the first paragraph: indentation and
#    paragraph border (space-only border, matching the style already used
#    by the third paragraph in this document).
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.LeftIndent = 11.25  # 225 twips = 11.25 pt

$b1 = $p1.Range.ParagraphFormat.Borders
$b1.DistanceFromTop = 5
$b1.DistanceFromLeft = 5
$b1.DistanceFromBottom = 5
$b1.DistanceFromRight = 5
